# Auto-generated edit script: refresh crypto price/volume snapshot
# (cryptos list updated by the scheduled GitHub Actions scraper).
# Row 12/13 also swap (WrappedliquidstakedEther2.0 <-> Polkadot reorder).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay literal text (e.g. trailing zeros,
# thousands-dot formatting like "27.073.12") so force Text format before
# writing, otherwise Excel would auto-coerce them to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.940.88'
$ws.Range("E2").Value = '  -1.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.693.00'
$ws.Range("E3").Value = '  -0.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.37%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.26'
$ws.Range("E5").Value = '  -0.58%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5249'
$ws.Range("E6").Value = '  -0.97%  '

# Row 7
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06602'
$ws.Range("E8").Value = '  +0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2613'
$ws.Range("E9").Value = '  -1.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.45'
$ws.Range("E10").Value = '  -2.03%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07705'
$ws.Range("E11").Value = '  +0.69%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.404'
$ws.Range("E12").Value = '  -3.83%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.929.93'
$ws.Range("E13").Value = '  -0.67%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.701.42'
$ws.Range("E14").Value = '  -0.38%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5683'
$ws.Range("E15").Value = '  -1.09%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8109'
$ws.Range("E16").Value = '  -0.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.15'
$ws.Range("E17").Value = '  -1.87%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.009.84'
$ws.Range("E18").Value = '  -0.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.09'
$ws.Range("E19").Value = '  +0.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.568'
$ws.Range("E21").Value = '  -2.31%  '

# Row 22
$ws.Range("E22").Value = '  -1.94%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.984'
$ws.Range("E23").Value = '  +0.26%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.007'
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.90'
$ws.Range("E25").Value = '  +1.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.726'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1190'
$ws.Range("E27").Value = '  -2.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.120'
$ws.Range("E28").Value = '  -2.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.99'
$ws.Range("E29").Value = '  -1.96%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05306'
$ws.Range("E30").Value = '  -1.35%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  -0.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.418'
$ws.Range("E32").Value = '  -2.36%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.301'
$ws.Range("E33").Value = '  -3.49%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.621'
$ws.Range("E34").Value = '  -1.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.799'
$ws.Range("E35").Value = '  -2.85%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9406'
$ws.Range("E36").Value = '  -1.13%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.393'
$ws.Range("E37").Value = '  -1.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5837'
$ws.Range("E38").Value = '  -0.60%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.174.11'
$ws.Range("E39").Value = '  +12.65%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01629'
$ws.Range("E40").Value = '  +0.16%  '

# Row 41
$ws.Range("E41").Value = '  +0.20%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.730'
$ws.Range("E42").Value = '  -2.67%  '

# Row 43
$ws.Range("E43").Value = '  -0.44%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.36'
$ws.Range("E44").Value = '  -0.81%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.837.95'
$ws.Range("E45").Value = '  -0.65%  '

# Row 46
$ws.Range("E46").Value = '  -2.73%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.78'
$ws.Range("E47").Value = '  -2.20%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4514'
$ws.Range("E48").Value = '  +0.48%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.005'
$ws.Range("E49").Value = '  +0.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.076'
$ws.Range("E50").Value = '  -0.27%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05196'
$ws.Range("E51").Value = '  -0.79%  '
